$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L4").Value = 317.95    # ALVARADO BAEZ EVELYN MARY - PIEDRA SINTERIZADA
$ws1.Range("M42").Value = 1438.28  # QUIROZ PEÑAFIEL DIANA ISABEL - PORCELANATO
$ws1.Range("M49").Value = 1582.59  # SANCHEZ CORREA MARCO EDUARDO - PORCELANATO
$ws1.Range("L52").Value = 0        # SISA GUANO CARLOS ALBERTO - PIEDRA SINTERIZADA

# --- Sheet "VENTA MENSUAL" ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F4").Value = 317.95    # ALVARADO BAEZ EVELYN MARY - noviembre
$ws2.Range("F42").Value = 1438.28  # QUIROZ PEÑAFIEL DIANA ISABEL - noviembre
$ws2.Range("F49").Value = 1582.59  # SANCHEZ CORREA MARCO EDUARDO - noviembre
$ws2.Range("F52").Value = 0        # SISA GUANO CARLOS ALBERTO - noviembre
$ws2.Range("F60").Value = 11829.48 # TOTAL - noviembre

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D11").Value = 2729.65             # PIEDRA SINTERIZADA - VENTA
$ws3.Range("E11").Value = 13418.35            # PIEDRA SINTERIZADA - POR CUMPLIR
$ws3.Range("F11").Value = 0.1690395095367848  # PIEDRA SINTERIZADA - CUMPLIMIENTO

$ws3.Range("D12").Value = 5690.91             # PORCELANATO - VENTA
$ws3.Range("E12").Value = 44616.09            # PORCELANATO - POR CUMPLIR
$ws3.Range("F12").Value = 0.113123620967261   # PORCELANATO - CUMPLIMIENTO

$ws3.Range("D14").Value = 13446.79            # TOTAL - VENTA
$ws3.Range("E14").Value = 84415.09766749099   # TOTAL - POR CUMPLIR
$ws3.Range("F14").Value = 0.1374057901446645  # TOTAL - CUMPLIMIENTO

$wb.Save()
